$d = $word.ActiveDocument

$replacements = @(
    @("2024-01-26 Friday", "2024-01-27 Saturday"),
    @("312×2=624", "544×8=4352"),
    @("977×7=6839", "753×2=1506"),
    @("732×9=6588", "768×7=5376"),
    @("547×2=1094", "165×7=1155"),
    @("661×6=3966", "713×3=2139"),
    @("723×3=2169", "134×8=1072"),
    @("967×5=4835", "364×5=1820"),
    @("849×2=1698", "970×9=8730"),
    @("128×8=1024", "413×3=1239"),
    @("865×6=5190", "706×7=4942"),
    @("882×4=3528", "609×9=5481"),
    @("518×3=1554", "815×3=2445"),
    @("359×5=1795", "465×3=1395"),
    @("915×7=6405", "742×3=2226"),
    @("442×3=1326", "888×7=6216"),
    @("712×3=2136", "306×6=1836"),
    @("813×9=7317", "111×4=444"),
    @("320×7=2240", "708×6=4248"),
    @("853×9=7677", "472×3=1416"),
    @("359×6=2154", "788×2=1576"),
    @("986×7=6902", "175×3=525"),
    @("655×4=2620", "333×4=1332"),
    @("517×5=2585", "585×8=4680"),
    @("956×6=5736", "136×2=272"),
    @("748×8=5984", "329×3=987")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
